$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 712.9167
$ws.Range("J17").Value = 718.1429000000001
$ws.Range("L17").Value = 2154.4287
$ws.Range("N17").Value = -2490.4287
$ws.Range("H46").Value = 5725
$ws.Range("I46").Value = 187.5
$ws.Range("K46").Value = 562.5
$ws.Range("M46").Value = -443.5
$ws.Range("H49").Value = 597.1429000000001
$ws.Range("I49").Value = 360
$ws.Range("J49").Value = 775
$ws.Range("K49").Value = 1080
$ws.Range("L49").Value = 2325
$ws.Range("M49").Value = -944
$ws.Range("N49").Value = -2597
$ws.Range("H60").Value = 5725
$ws.Range("I60").Value = 187.5
$ws.Range("K60").Value = 562.5
$ws.Range("M60").Value = -78.5
$ws.Range("H63").Value = 45000
$ws.Range("J63").Value = 45000
$ws.Range("L63").Value = 45000
$ws.Range("N63").Value = -46248
$ws.Range("H66").Value = 45000
$ws.Range("J66").Value = 45000
$ws.Range("L66").Value = 135000
$ws.Range("N66").Value = -141240
$ws.Range("H106").Value = 50100844
$ws.Range("I106").Value = 112047.22
$ws.Range("K106").Value = 112047.22
$ws.Range("M106").Value = -111416.22
$ws.Range("H112").Value = 1540.7142
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 1540.7142
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 4622.142599999999
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -6838.142599999999
$ws.Range("H132").Value = 2013526
$ws.Range("I132").Value = 2464578.8
$ws.Range("J132").Value = 1136.5385
$ws.Range("K132").Value = 7393736.399999999
$ws.Range("L132").Value = 3409.6155
$ws.Range("M132").Value = -7391206.399999999
$ws.Range("N132").Value = -8469.6155
$ws.Range("H137").Value = 1289.175
$ws.Range("I137").Value = 941.8333
$ws.Range("J137").Value = 2331.2
$ws.Range("K137").Value = 2825.4999
$ws.Range("L137").Value = 6993.599999999999
$ws.Range("M137").Value = -275.4998999999998
$ws.Range("N137").Value = -12093.6
$ws.Range("H138").Value = 1868.9692
$ws.Range("I138").Value = 1594.9667
$ws.Range("J138").Value = 2103.8286
$ws.Range("K138").Value = 4784.9001
$ws.Range("L138").Value = 6311.485799999999
$ws.Range("M138").Value = 355.0999000000002
$ws.Range("N138").Value = -16591.4858
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1264.5555
$ws.Range("I45").Value = 1201.8334
$ws.Range("J45").Value = 1390
$ws.Range("K45").Value = 1201.8334
$ws.Range("L45").Value = 1390
$ws.Range("M45").Value = -824.8334
$ws.Range("N45").Value = -2144
$ws.Range("H61").Value = 1591.5581
$ws.Range("I61").Value = 1161.3334
$ws.Range("J61").Value = 3011.3
$ws.Range("K61").Value = 1161.3334
$ws.Range("L61").Value = 3011.3
$ws.Range("M61").Value = -949.3334
$ws.Range("N61").Value = -3435.3
$ws.Range("H122").Value = 1313.4
$ws.Range("I122").Value = 1313.4
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3940.2
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1490.2
$ws.Range("N122").ClearContents()
$ws.Range("H136").Value = 1591.5581
$ws.Range("I136").Value = 1161.3334
$ws.Range("J136").Value = 3011.3
$ws.Range("K136").Value = 3484.0002
$ws.Range("L136").Value = 9033.900000000001
$ws.Range("M136").Value = -934.0001999999999
$ws.Range("N136").Value = -14133.9
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1352.9231
$ws.Range("I94").Value = 807.25
$ws.Range("K94").Value = 807.25
$ws.Range("M94").Value = -356.25
$ws.Range("H107").Value = 1203.6666
$ws.Range("I107").Value = 1205.5
$ws.Range("J107").Value = 1200
$ws.Range("K107").Value = 1205.5
$ws.Range("L107").Value = 1200
$ws.Range("M107").Value = 714.5
$ws.Range("N107").Value = -5040
$ws.Range("H134").Value = 22994.17
$ws.Range("I134").Value = 31865.484
$ws.Range("J134").Value = 2083.2144
$ws.Range("K134").Value = 95596.452
$ws.Range("L134").Value = 6249.6432
$ws.Range("M134").Value = -93061.452
$ws.Range("N134").Value = -11319.6432
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 20204556
$ws.Range("I31").Value = 23811182
$ws.Range("J31").Value = 13892961
$ws.Range("K31").Value = 23811182
$ws.Range("L31").Value = 13892961
$ws.Range("M31").Value = -23810887
$ws.Range("N31").Value = -13893551
$ws.Range("H34").Value = 20204556
$ws.Range("I34").Value = 23811182
$ws.Range("J34").Value = 13892961
$ws.Range("K34").Value = 23811182
$ws.Range("L34").Value = 13892961
$ws.Range("M34").Value = -23810980
$ws.Range("N34").Value = -13893365
$ws.Range("H64").Value = 40000
$ws.Range("J64").Value = 40000
$ws.Range("L64").Value = 40000
$ws.Range("N64").Value = -40496
$ws.Range("H67").Value = 40000
$ws.Range("J67").Value = 40000
$ws.Range("L67").Value = 40000
$ws.Range("N67").Value = -41716
$ws.Range("H132").Value = 4220643.5
$ws.Range("I132").Value = 6061803.5
$ws.Range("J132").Value = 1318.5416
$ws.Range("K132").Value = 18185410.5
$ws.Range("L132").Value = 3955.6248
$ws.Range("M132").Value = -18182880.5
$ws.Range("N132").Value = -9015.6248
$ws.Range("H134").Value = 831.0925999999999
$ws.Range("I134").Value = 789.5098
$ws.Range("J134").Value = 1538
$ws.Range("K134").Value = 2368.5294
$ws.Range("L134").Value = 4614
$ws.Range("M134").Value = 166.4705999999996
$ws.Range("N134").Value = -9684
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 55556300
$ws.Range("I129").Value = 1115
$ws.Range("K129").Value = 3345
$ws.Range("M129").Value = 1655
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3555.75
$ws.Range("I80").Value = 2394.625
$ws.Range("J80").Value = 5878
$ws.Range("K80").Value = 2394.625
$ws.Range("L80").Value = 5878
$ws.Range("M80").Value = -1396.625
$ws.Range("N80").Value = -7874
$ws.Range("H83").Value = 3555.75
$ws.Range("I83").Value = 2394.625
$ws.Range("J83").Value = 5878
$ws.Range("K83").Value = 11973.125
$ws.Range("L83").Value = 29390
$ws.Range("M83").Value = -6981.125
$ws.Range("N83").Value = -39374
$ws.Range("H97").Value = 1040
$ws.Range("I97").Value = 890
$ws.Range("K97").Value = 890
$ws.Range("M97").Value = -394
$ws.Range("H122").Value = 40003560
$ws.Range("I122").Value = 76927190
$ws.Range("J122").Value = 2954.8333
$ws.Range("K122").Value = 230781570
$ws.Range("L122").Value = 8864.499899999999
$ws.Range("M122").Value = -230779120
$ws.Range("N122").Value = -13764.4999
$ws.Range("H132").Value = 22986.447
$ws.Range("I132").Value = 24750.303
$ws.Range("J132").Value = 4025
$ws.Range("K132").Value = 74250.909
$ws.Range("L132").Value = 12075
$ws.Range("M132").Value = -71720.909
$ws.Range("N132").Value = -17135
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 711.125
$ws.Range("I16").Value = 363.70587
$ws.Range("J16").Value = 1554.8572
$ws.Range("K16").Value = 363.70587
$ws.Range("L16").Value = 1554.8572
$ws.Range("M16").Value = -193.70587
$ws.Range("N16").Value = -1894.8572
$ws.Range("H82").Value = 1366
$ws.Range("I82").Value = 1395.3636
$ws.Range("J82").Value = 1312.1666
$ws.Range("K82").Value = 1395.3636
$ws.Range("L82").Value = 1312.1666
$ws.Range("M82").Value = -1034.3636
$ws.Range("N82").Value = -2034.1666
$ws.Range("H85").Value = 1366
$ws.Range("I85").Value = 1395.3636
$ws.Range("J85").Value = 1312.1666
$ws.Range("K85").Value = 1395.3636
$ws.Range("L85").Value = 1312.1666
$ws.Range("M85").Value = -147.3635999999999
$ws.Range("N85").Value = -3808.1666
$ws.Range("H132").Value = 8339891.5
$ws.Range("I132").Value = 12829548
$ws.Range("J132").Value = 1956.7142
$ws.Range("K132").Value = 38488644
$ws.Range("L132").Value = 5870.142599999999
$ws.Range("M132").Value = -38486114
$ws.Range("N132").Value = -10930.1426
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3056.8572
$ws.Range("I81").Value = 4049.75
$ws.Range("J81").Value = 1733
$ws.Range("K81").Value = 8099.5
$ws.Range("L81").Value = 3466
$ws.Range("M81").Value = -7038.5
$ws.Range("N81").Value = -5588
$ws.Range("H84").Value = 3056.8572
$ws.Range("I84").Value = 4049.75
$ws.Range("J84").Value = 1733
$ws.Range("K84").Value = 40497.5
$ws.Range("L84").Value = 17330
$ws.Range("M84").Value = -35193.5
$ws.Range("N84").Value = -27938
$ws.Range("H96").Value = 50000760
$ws.Range("I96").Value = 62500700
$ws.Range("J96").Value = 1004
$ws.Range("K96").Value = 62500700
$ws.Range("L96").Value = 1004
$ws.Range("M96").Value = -62499327
$ws.Range("N96").Value = -3750
$ws.Range("H126").Value = 7275.8823
$ws.Range("I126").Value = 12951.333
$ws.Range("K126").Value = 38853.999
$ws.Range("M126").Value = -36383.999
